# Generate Report for Handback
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Column-width constants. Excel's ColumnWidth setter quantizes the request
# to the nearest 1/6 character unit and then stores (value + 5/6) as the
# OOXML `width`; feed it (target - 5/6) so the stored width lands on the
# requested value.
$wideColWidth  = 29.144371396019366   # -> stored width 29.9777047293527 (closest snap: 30)
$fortyColWidth = 39.166666666666664   # -> stored width 40

# --- Overview sheet: widen columns E and F ---
$wsOverview.Columns.Item(5).ColumnWidth = $wideColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $wideColWidth

# --- zh-cn sheet ---
# Status text update (shared string also used by de-de sheet)
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"

# Widen columns C, I, J
$wsZhCn.Columns.Item(3).ColumnWidth = $wideColWidth
$wsZhCn.Columns.Item(9).ColumnWidth = $fortyColWidth
$wsZhCn.Columns.Item(10).ColumnWidth = $fortyColWidth

# Latest Target File (I2) now points to the handback file, styled as a hyperlink
$wsZhCn.Range("I2").Value = "af782309-5a61-44cf-acfc-13ae29cb091c.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7b799d87378b4ed77b92f81499b8589a8383db71/e2e/af782309-5a61-44cf-acfc-13ae29cb091c.md", [System.Type]::Missing, "af782309-5a61-44cf-acfc-13ae29cb091c.md", "af782309-5a61-44cf-acfc-13ae29cb091c.md") | Out-Null

# Latest Handback File (J2)
$wsZhCn.Range("J2").Value = "af782309-5a61-44cf-acfc-13ae29cb091c.684a6c4b50f9301899e07b0a7b40f9f86aea1a87.zh-cn.xlf"

# Latest Handback DateTime (K2)
$wsZhCn.Range("K2").Value = "2016-09-02 15:13:51"

# --- de-de sheet ---
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"

$wsDeDe.Columns.Item(3).ColumnWidth = $wideColWidth
$wsDeDe.Columns.Item(9).ColumnWidth = $fortyColWidth
$wsDeDe.Columns.Item(10).ColumnWidth = $fortyColWidth

$wsDeDe.Range("I2").Value = "af782309-5a61-44cf-acfc-13ae29cb091c.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7b799d87378b4ed77b92f81499b8589a8383db71/e2e/af782309-5a61-44cf-acfc-13ae29cb091c.md", [System.Type]::Missing, "af782309-5a61-44cf-acfc-13ae29cb091c.md", "af782309-5a61-44cf-acfc-13ae29cb091c.md") | Out-Null

$wsDeDe.Range("J2").Value = "af782309-5a61-44cf-acfc-13ae29cb091c.684a6c4b50f9301899e07b0a7b40f9f86aea1a87.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-09-02 15:13:58"
